$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 570 (the "「我が魂と共に」..." post), shifting all subsequent rows up by one.
$ws.Rows.Item(570).Delete()
